# Update exported bill amounts.
# The source cells hold their amounts as text (numeric-looking strings),
# so we force each target cell to a Text number format before writing the
# new value -- this keeps Excel from "smart" converting the literal
# "1175.10" into the number 1175.1, which would silently drop the
# trailing zero / change the cell's underlying type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "H2" = "1175.10"
    "I2" = "1175.10"
    "H3" = "1125.45"
    "I3" = "1125.45"
    "I4" = "2300.55"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
